$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$nl = [char]11

$tbl.Cell(1, 1).Range.Text = "21 x 11" + $nl + "  1    1" + $nl + "  ----" + $nl + "2|    |" + $nl + "1|    |"
$tbl.Cell(1, 2).Range.Text = "98 x 57" + $nl + "  5    7" + $nl + "  ----" + $nl + "9|    |" + $nl + "8|    |"
$tbl.Cell(1, 3).Range.Text = "29 x 68" + $nl + "  6    8" + $nl + "  ----" + $nl + "2|    |" + $nl + "9|    |"
$tbl.Cell(2, 1).Range.Text = "77 x 61" + $nl + "  6    1" + $nl + "  ----" + $nl + "7|    |" + $nl + "7|    |"
$tbl.Cell(2, 2).Range.Text = "81 x 41" + $nl + "  4    1" + $nl + "  ----" + $nl + "8|    |" + $nl + "1|    |"
$tbl.Cell(2, 3).Range.Text = "90 x 28" + $nl + "  2    8" + $nl + "  ----" + $nl + "9|    |" + $nl + "0|    |"
$tbl.Cell(3, 1).Range.Text = "82 x 34" + $nl + "  3    4" + $nl + "  ----" + $nl + "8|    |" + $nl + "2|    |"
$tbl.Cell(3, 2).Range.Text = "65 x 38" + $nl + "  3    8" + $nl + "  ----" + $nl + "6|    |" + $nl + "5|    |"
$tbl.Cell(3, 3).Range.Text = "29 x 90" + $nl + "  9    0" + $nl + "  ----" + $nl + "2|    |" + $nl + "9|    |"
$tbl.Cell(4, 1).Range.Text = "67 x 83" + $nl + "  8    3" + $nl + "  ----" + $nl + "6|    |" + $nl + "7|    |"
$tbl.Cell(4, 2).Range.Text = "69 x 10" + $nl + "  1    0" + $nl + "  ----" + $nl + "6|    |" + $nl + "9|    |"
$tbl.Cell(4, 3).Range.Text = "93 x 70" + $nl + "  7    0" + $nl + "  ----" + $nl + "9|    |" + $nl + "3|    |"
$tbl.Cell(5, 1).Range.Text = "43 x 52" + $nl + "  5    2" + $nl + "  ----" + $nl + "4|    |" + $nl + "3|    |"
$tbl.Cell(5, 2).Range.Text = "90 x 64" + $nl + "  6    4" + $nl + "  ----" + $nl + "9|    |" + $nl + "0|    |"
$tbl.Cell(5, 3).Range.Text = "75 x 93" + $nl + "  9    3" + $nl + "  ----" + $nl + "7|    |" + $nl + "5|    |"
